$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '65.707.49'
$ws.Range('E2').Value = '  +0.06%  '
# Row 3
$ws.Range('D3').Value = '2.660.12'
$ws.Range('E3').Value = '  -0.38%  '
# Row 4
$ws.Range('E4').Value = '  +0.10%  '
# Row 5
$ws.Range('D5').Value = "'596.41"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.00%  '
# Row 6
$ws.Range('D6').Value = "'157.17"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.06%  '
# Row 7
$ws.Range('E7').Value = '  +4.53%  '
# Row 8
$ws.Range('E8').Value = '  +0.04%  '
# Row 9
$ws.Range('E9').Value = '  -4.26%  '
# Row 10
$ws.Range('D10').Value = "'0.398"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.92%  '
# Row 11
$ws.Range('E11').Value = '  -0.58%  '
# Row 12
$ws.Range('E12').Value = '  +1.29%  '
# Row 13
$ws.Range('D13').Value = "'28.63"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.99%  '
# Row 14
$ws.Range('D14').Value = "'0.0000195"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.34%  '
# Row 15
$ws.Range('D15').Value = '3.138.22'
$ws.Range('E15').Value = '  -0.41%  '
# Row 16
$ws.Range('D16').Value = '65.578.43'
$ws.Range('E16').Value = '  +0.08%  '
# Row 17
$ws.Range('D17').Value = '2.685.15'
$ws.Range('E17').Value = '  +1.90%  '
# Row 18
$ws.Range('D18').Value = "'12.55"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.91%  '
# Row 19
$ws.Range('D19').Value = "'4.77"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.78%  '
# Row 20
$ws.Range('D20').Value = "'349.64"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.71%  '
# Row 21
$ws.Range('D21').Value = "'7.44"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.46%  '
# Row 22
$ws.Range('E22').Value = '  +0.16%  '
# Row 23
$ws.Range('D23').Value = "'69.62"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.30%  '
# Row 24
$ws.Range('D24').Value = "'1.79"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +8.34%  '
# Row 25
$ws.Range('E25').Value = '  +0.81%  '
# Row 26
$ws.Range('D26').Value = "'9.52"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.55%  '
# Row 27
$ws.Range('E27').Value = '  +1.86%  '
# Row 28
$ws.Range('D28').Value = "'560.76"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.60%  '
# Row 29
$ws.Range('D29').Value = "'8.04"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.67%  '
# Row 30
$ws.Range('D30').Value = "'0.163"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.70%  '
# Row 31
$ws.Range('E31').Value = '  -0.15%  '
# Row 32
$ws.Range('D32').Value = "'2.13"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.04%  '
# Row 33
$ws.Range('D33').Value = "'1.81"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.68%  '
# Row 34
$ws.Range('D34').Value = "'6.54"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.09%  '
# Row 35
$ws.Range('E35').Value = '  -0.84%  '
# Row 36
$ws.Range('E36').Value = '  -1.20%  '
# Row 37
$ws.Range('D37').Value = "'20.43"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.27%  '
# Row 38
$ws.Range('D38').Value = "'1.00"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.14%  '
# Row 39
$ws.Range('E39').Value = '  -0.13%  '
# Row 40
$ws.Range('D40').Value = "'155.77"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.06%  '
# Row 41
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').Value = "'160.68"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.74%  '
# Row 42
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').Value = "'4.06"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.32%  '
# Row 43
$ws.Range('B43').Value = 'Hedera'
$ws.Range('C43').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D43').Value = "'0.0605"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.90%  '
# Row 44
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').Value = "'2.27"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.80%  '
# Row 45
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').Value = "'22.59"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.36%  '
# Row 46
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').Value = "'0.640"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.60%  '
# Row 47
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').Value = "'0.0254"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.94%  '
# Row 48
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').Value = "'0.102"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.27%  '
# Row 49
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = "'19.75"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.26%  '
# Row 50
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0245'
$ws.Range('E50').Value = '  +3.09%  '
# Row 51
$ws.Range('B51').Value = 'ONDO'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D51').Value = "'0.802"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.74%  '
